$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2750
$ws.Range("J40").Value = 2750
$ws.Range("L40").Value = 2750
$ws.Range("N40").Value = -3100
$ws.Range("H76").Value = 7672357.5
$ws.Range("I76").Value = 4548.5
$ws.Range("K76").Value = 4548.5
$ws.Range("M76").Value = -4233.5
$ws.Range("H79").Value = 7672357.5
$ws.Range("I79").Value = 4548.5
$ws.Range("K79").Value = 4548.5
$ws.Range("M79").Value = -3456.5
$ws.Range("H80").Value = 994740.4399999999
$ws.Range("J80").Value = 3002.9092
$ws.Range("L80").Value = 9008.7276
$ws.Range("N80").Value = -11004.7276
$ws.Range("H83").Value = 994740.4399999999
$ws.Range("J83").Value = 3002.9092
$ws.Range("L83").Value = 27026.1828
$ws.Range("N83").Value = -37010.1828
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H105").Value = 28845.375
$ws.Range("J105").Value = 28845.375
$ws.Range("L105").Value = 28845.375
$ws.Range("N105").Value = -35833.375
$ws.Range("H112").Value = 1668.7059
$ws.Range("I112").Value = 694.8
$ws.Range("K112").Value = 2084.4
$ws.Range("M112").Value = -976.3999999999996
$ws.Range("H138").Value = 2434.3044
$ws.Range("I138").Value = 1997.375
$ws.Range("J138").Value = 2667.3333
$ws.Range("K138").Value = 5992.125
$ws.Range("L138").Value = 8001.999899999999
$ws.Range("M138").Value = -852.125
$ws.Range("N138").Value = -18281.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 995
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H32").Value = 7173.3696
$ws.Range("I32").Value = 3345
$ws.Range("J32").Value = 14989.625
$ws.Range("K32").Value = 3345
$ws.Range("L32").Value = 14989.625
$ws.Range("M32").Value = -3058
$ws.Range("N32").Value = -15563.625
$ws.Range("H45").Value = 2397
$ws.Range("I45").Value = 2496.25
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2496.25
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -2119.25
$ws.Range("N45").Value = -2754
$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H116").Value = 995
$ws.Range("I116").Value = 1000
$ws.Range("K116").Value = 1000
$ws.Range("M116").Value = 1294
$ws.Range("H135").Value = 103599.6
$ws.Range("J135").Value = 103599.6
$ws.Range("L135").Value = 103599.6
$ws.Range("N135").Value = -113739.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 995
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -886
$ws.Range("H109").Value = 39894
$ws.Range("J109").Value = 39894
$ws.Range("L109").Value = 39894
$ws.Range("N109").Value = -42668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5041.9546
$ws.Range("I31").Value = 13914.2
$ws.Range("J31").Value = 3457.625
$ws.Range("K31").Value = 13914.2
$ws.Range("L31").Value = 3457.625
$ws.Range("M31").Value = -13619.2
$ws.Range("N31").Value = -4047.625
$ws.Range("H34").Value = 5041.9546
$ws.Range("I34").Value = 13914.2
$ws.Range("J34").Value = 3457.625
$ws.Range("K34").Value = 13914.2
$ws.Range("L34").Value = 3457.625
$ws.Range("M34").Value = -13712.2
$ws.Range("N34").Value = -3861.625
$ws.Range("H43").Value = 12378
$ws.Range("J43").Value = 12378
$ws.Range("L43").Value = 12378
$ws.Range("N43").Value = -12746
$ws.Range("H55").Value = 10020
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9685
$ws.Range("H101").Value = 12378
$ws.Range("J101").Value = 12378
$ws.Range("L101").Value = 12378
$ws.Range("N101").Value = -18868
$ws.Range("H105").Value = 852.8333
$ws.Range("I105").Value = 700
$ws.Range("J105").Value = 883.4
$ws.Range("K105").Value = 700
$ws.Range("L105").Value = 883.4
$ws.Range("M105").Value = 1047
$ws.Range("N105").Value = -4377.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16532376
$ws.Range("I4").Value = 37652330
$ws.Range("K4").Value = 112956990
$ws.Range("M4").Value = -112956878
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 600
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -938
$ws.Range("H34").Value = 1691.8462
$ws.Range("J34").Value = 2643.5
$ws.Range("L34").Value = 7930.5
$ws.Range("N34").Value = -8098.5
$ws.Range("H39").Value = 7772.4546
$ws.Range("J39").Value = 4249.625
$ws.Range("L39").Value = 12748.875
$ws.Range("N39").Value = -13336.875
$ws.Range("H55").Value = 400
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4600.3
$ws.Range("I80").Value = 4199.6
$ws.Range("J80").Value = 5001
$ws.Range("K80").Value = 4199.6
$ws.Range("L80").Value = 5001
$ws.Range("M80").Value = -3201.6
$ws.Range("N80").Value = -6997
$ws.Range("H83").Value = 4600.3
$ws.Range("I83").Value = 4199.6
$ws.Range("J83").Value = 5001
$ws.Range("K83").Value = 20998
$ws.Range("L83").Value = 25005
$ws.Range("M83").Value = -16006
$ws.Range("N83").Value = -34989

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 894.1818
$ws.Range("I22").Value = 859.4
$ws.Range("J22").Value = 923.1667
$ws.Range("K22").Value = 859.4
$ws.Range("L22").Value = 923.1667
$ws.Range("M22").Value = -564.4
$ws.Range("N22").Value = -1513.1667
$ws.Range("H27").Value = 894.1818
$ws.Range("I27").Value = 859.4
$ws.Range("J27").Value = 923.1667
$ws.Range("K27").Value = 859.4
$ws.Range("L27").Value = 923.1667
$ws.Range("M27").Value = -752.4
$ws.Range("N27").Value = -1137.1667
$ws.Range("H48").Value = 34999.5
$ws.Range("I48").Value = 29999
$ws.Range("J48").Value = 40000
$ws.Range("K48").Value = 29999
$ws.Range("L48").Value = 40000
$ws.Range("M48").Value = -29338
$ws.Range("N48").Value = -41322
$ws.Range("H68").Value = 14184.875
$ws.Range("I68").Value = 5596
$ws.Range("J68").Value = 28499.666
$ws.Range("K68").Value = 5596
$ws.Range("L68").Value = 28499.666
$ws.Range("M68").Value = -4847
$ws.Range("N68").Value = -29997.666
$ws.Range("H71").Value = 14184.875
$ws.Range("I71").Value = 5596
$ws.Range("J71").Value = 28499.666
$ws.Range("K71").Value = 27980
$ws.Range("L71").Value = 142498.33
$ws.Range("M71").Value = -24236
$ws.Range("N71").Value = -149986.33
$ws.Range("H74").Value = 63951
$ws.Range("I74").Value = 49938.75
$ws.Range("K74").Value = 49938.75
$ws.Range("M74").Value = -48940.75
$ws.Range("H77").Value = 63951
$ws.Range("I77").Value = 49938.75
$ws.Range("K77").Value = 149816.25
$ws.Range("M77").Value = -144824.25
$ws.Range("H103").Value = 18968.2
$ws.Range("J103").Value = 18968.2
$ws.Range("L103").Value = 18968.2
$ws.Range("N103").Value = -21312.2
$ws.Range("H110").Value = 61153.668
$ws.Range("J110").Value = 61153.668
$ws.Range("L110").Value = 61153.668
$ws.Range("N110").Value = -69333.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2171139
$ws.Range("I62").Value = 11907262
$ws.Range("J62").Value = 7556.1113
$ws.Range("K62").Value = 11907262
$ws.Range("L62").Value = 7556.1113
$ws.Range("M62").Value = -11906638
$ws.Range("N62").Value = -8804.1113
$ws.Range("H65").Value = 2171139
$ws.Range("I65").Value = 11907262
$ws.Range("J65").Value = 7556.1113
$ws.Range("K65").Value = 59536310
$ws.Range("L65").Value = 37780.5565
$ws.Range("M65").Value = -59533190
$ws.Range("N65").Value = -44020.5565
$ws.Range("H104").Value = 9604.833000000001
$ws.Range("J104").Value = 9604.833000000001
$ws.Range("L104").Value = 9604.833000000001
$ws.Range("N104").Value = -16592.833
